$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the corrupted header cell: A1 held the whole recorded Python
# script (copy/paste mistake) glued onto the real header text "Friend1".
# Restore it to the plain header label.
$ws.Range("A1").Value = "Friend1"

# --- Fix two bad lookups in the friend-pair table: B5 and B14 pointed at
# stray "Jack"/"Paul" entries instead of "Emily".
$ws.Range("B5").Value = "Emily"
$ws.Range("B14").Value = "Emily"

# --- Formatting bug fix: header row + the A:B column default alignment
# was left as "General"; make it explicitly "Left" (matches the rest of
# the data rows, which already use Left).
$ws.Range("A1:B1").HorizontalAlignment = -4131
$ws.Columns("A:B").HorizontalAlignment = -4131

Write-Output "done"
